# Update countries & provincias Spain
# Refresh the COVID "Pais" table: a handful of countries (Bolivia/China,
# Venezuela/Austria, Timor Oriental/Santa Lucia) swap rank as their totals
# are refreshed, a few other rows get updated figures, and the "last
# updated" timestamp moves from 03:20 to 04:37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Title / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 04:37"

# --- Rows 32/33: Bolivia overtakes China ------------------------------
$ws.Range("A32").Value = "Bolivia"
$ws.Range("B32").Value = 85141
$ws.Range("C32").Value = 1780
$ws.Range("D32").Value = 26437
$ws.Range("E32").Value = 55319
$ws.Range("G32").Value = 65
$ws.Range("H32").Value = 3385

$ws.Range("A33").Value = "China"
$ws.Range("B33").Value = 84528
$ws.Range("C33").Value = 37
$ws.Range("D33").Value = 79057
$ws.Range("E33").Value = 837
$ws.Range("H33").Value = 4634

# --- Rows 67/68: Venezuela overtakes Austria --------------------------
$ws.Range("A67").Value = "Venezuela"
$ws.Range("B67").Value = 22299
$ws.Range("D67").Value = 12146
$ws.Range("E67").Value = 9958
$ws.Range("H67").Value = 195

$ws.Range("A68").Value = "Austria"
$ws.Range("B68").Value = 21566
$ws.Range("D68").Value = 19464
$ws.Range("E68").Value = 1383
$ws.Range("H68").Value = 719

# --- Row 72: Australia figures refreshed (no rank change) -------------
$ws.Range("B72").Value = 19890
$ws.Range("C72").Value = 446
$ws.Range("D72").Value = 10941
$ws.Range("E72").Value = 8694
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 255

# --- Row 77: Corea del Sur figures refreshed (no rank change) ---------
$ws.Range("B77").Value = 14499
$ws.Range("C77").Value = 43
$ws.Range("D77").Value = 13501
$ws.Range("E77").Value = 696

# --- Row 138: Nueva Zelanda minor correction ---------------------------
$ws.Range("D138").Value = 1524
$ws.Range("E138").Value = 23

# --- Rows 202/203: Timor Oriental overtakes Santa Lucia ----------------
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
